$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 524117.66
$ws.Range("I88").Value = 1120539.6
$ws.Range("J88").Value = 2248.5
$ws.Range("K88").Value = 1120539.6
$ws.Range("L88").Value = 2248.5
$ws.Range("M88").Value = -1120133.6
$ws.Range("N88").Value = -3060.5
$ws.Range("H91").Value = 524117.66
$ws.Range("I91").Value = 1120539.6
$ws.Range("J91").Value = 2248.5
$ws.Range("K91").Value = 1120539.6
$ws.Range("L91").Value = 2248.5
$ws.Range("M91").Value = -1119135.6
$ws.Range("N91").Value = -5056.5
$ws.Range("H92").Value = 471.5625
$ws.Range("I92").Value = 483
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 483
$ws.Range("L92").Value = 300
$ws.Range("M92").Value = 765
$ws.Range("N92").Value = -2796
$ws.Range("H132").Value = 11500.549
$ws.Range("I132").Value = 1650.4937
$ws.Range("J132").Value = 67083
$ws.Range("K132").Value = 4951.4811
$ws.Range("L132").Value = 201249
$ws.Range("M132").Value = -2421.4811
$ws.Range("N132").Value = -206309
$ws.Range("H135").Value = 7813463
$ws.Range("I135").Value = 644.3488
$ws.Range("J135").Value = 23811138
$ws.Range("K135").Value = 5799.1392
$ws.Range("L135").Value = 214300242
$ws.Range("M135").Value = -3264.1392
$ws.Range("N135").Value = -214305312
$ws.Range("H137").Value = 2658.743
$ws.Range("I137").Value = 755.549
$ws.Range("J137").Value = 7767.316
$ws.Range("K137").Value = 2266.647
$ws.Range("L137").Value = 23301.948
$ws.Range("M137").Value = 283.3530000000001
$ws.Range("N137").Value = -28401.948
$ws.Range("H138").Value = 3589.6338
$ws.Range("I138").Value = 2168
$ws.Range("J138").Value = 4362.2607
$ws.Range("K138").Value = 6504
$ws.Range("L138").Value = 13086.7821
$ws.Range("M138").Value = -1364
$ws.Range("N138").Value = -23366.7821
$ws.Range("H141").Value = 833.42255
$ws.Range("I141").Value = 593.6667
$ws.Range("J141").Value = 9105
$ws.Range("K141").Value = 1781.0001
$ws.Range("L141").Value = 27315
$ws.Range("M141").Value = 3398.9999
$ws.Range("N141").Value = -37675

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2951.35
$ws.Range("I45").Value = 2392.6667
$ws.Range("J45").Value = 3408.4546
$ws.Range("K45").Value = 2392.6667
$ws.Range("L45").Value = 3408.4546
$ws.Range("M45").Value = -2015.6667
$ws.Range("N45").Value = -4162.4546
$ws.Range("H61").Value = 1194.57
$ws.Range("I61").Value = 1027.8182
$ws.Range("J61").Value = 1602.1852
$ws.Range("K61").Value = 1027.8182
$ws.Range("L61").Value = 1602.1852
$ws.Range("M61").Value = -815.8181999999999
$ws.Range("N61").Value = -2026.1852
$ws.Range("H74").Value = 1755.7727
$ws.Range("I74").Value = 1618.5
$ws.Range("K74").Value = 1618.5
$ws.Range("M74").Value = -744.5
$ws.Range("H77").Value = 1755.7727
$ws.Range("I77").Value = 1618.5
$ws.Range("K77").Value = 8092.5
$ws.Range("M77").Value = -3724.5
$ws.Range("H97").Value = 522.2222
$ws.Range("I97").Value = 409.43243
$ws.Range("K97").Value = 409.43243
$ws.Range("M97").Value = 86.56756999999999
$ws.Range("H110").Value = 1920.125
$ws.Range("I110").Value = 1920.125
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1920.125
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 124.875
$ws.Range("H132").Value = 6098959
$ws.Range("I132").Value = 8929525
$ws.Range("J132").Value = 2354.0386
$ws.Range("K132").Value = 26788575
$ws.Range("L132").Value = 7062.1158
$ws.Range("M132").Value = -26786045
$ws.Range("N132").Value = -12122.1158
$ws.Range("H136").Value = 1194.57
$ws.Range("I136").Value = 1027.8182
$ws.Range("J136").Value = 1602.1852
$ws.Range("K136").Value = 3083.4546
$ws.Range("L136").Value = 4806.5556
$ws.Range("M136").Value = -533.4546
$ws.Range("N136").Value = -9906.5556
$ws.Range("N110").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2304.3215
$ws.Range("I99").Value = 2068.6365
$ws.Range("J99").Value = 3168.5
$ws.Range("K99").Value = 2068.6365
$ws.Range("L99").Value = 3168.5
$ws.Range("M99").Value = -570.6365000000001
$ws.Range("N99").Value = -6164.5
$ws.Range("H105").Value = 2156.3684
$ws.Range("I105").Value = 2190
$ws.Range("K105").Value = 2190
$ws.Range("M105").Value = -443
$ws.Range("H134").Value = 2059.69
$ws.Range("I134").Value = 824.3148
$ws.Range("J134").Value = 3509.913
$ws.Range("K134").Value = 2472.9444
$ws.Range("L134").Value = 10529.739
$ws.Range("M134").Value = 62.05560000000014
$ws.Range("N134").Value = -15599.739

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H31").Value = 3335.27
$ws.Range("I31").Value = 1334
$ws.Range("J31").Value = 3745.1687
$ws.Range("K31").Value = 1334
$ws.Range("L31").Value = 3745.1687
$ws.Range("M31").Value = -1039
$ws.Range("N31").Value = -4335.1687
$ws.Range("H34").Value = 3335.27
$ws.Range("I34").Value = 1334
$ws.Range("J34").Value = 3745.1687
$ws.Range("K34").Value = 1334
$ws.Range("L34").Value = 3745.1687
$ws.Range("M34").Value = -1132
$ws.Range("N34").Value = -4149.1687
$ws.Range("H58").Value = 1017.2241
$ws.Range("I58").Value = 687.5833
$ws.Range("J58").Value = 2599.5
$ws.Range("K58").Value = 687.5833
$ws.Range("L58").Value = 2599.5
$ws.Range("M58").Value = -484.5833
$ws.Range("N58").Value = -3005.5
$ws.Range("H105").Value = 3542.5715
$ws.Range("I105").Value = 3616.4167
$ws.Range("J105").Value = 3099.5
$ws.Range("K105").Value = 3616.4167
$ws.Range("L105").Value = 3099.5
$ws.Range("M105").Value = -1869.4167
$ws.Range("N105").Value = -6593.5
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("H132").Value = 21314.285
$ws.Range("I132").Value = 931.5246
$ws.Range("K132").Value = 2794.5738
$ws.Range("M132").Value = -264.5738000000001
$ws.Range("H134").Value = 210105.83
$ws.Range("I134").Value = 705.9231
$ws.Range("J134").Value = 936025.5600000001
$ws.Range("K134").Value = 2117.7693
$ws.Range("L134").Value = 2808076.68
$ws.Range("M134").Value = 417.2307000000001
$ws.Range("N134").Value = -2813146.68
$ws.Range("H136").Value = 1017.2241
$ws.Range("I136").Value = 687.5833
$ws.Range("J136").Value = 2599.5
$ws.Range("K136").Value = 2062.7499
$ws.Range("L136").Value = 7798.5
$ws.Range("M136").Value = 487.2501000000002
$ws.Range("N136").Value = -12898.5
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4203.8184
$ws.Range("I5").Value = 8008.2856
$ws.Range("J5").Value = 1400.5264
$ws.Range("K5").Value = 24024.8568
$ws.Range("L5").Value = 4201.5792
$ws.Range("M5").Value = -23912.8568
$ws.Range("N5").Value = -4425.5792
$ws.Range("H80").Value = 143143730
$ws.Range("I80").Value = 1000151
$ws.Range("J80").Value = 200001150
$ws.Range("K80").Value = 3000453
$ws.Range("L80").Value = 600003450
$ws.Range("M80").Value = -2999517
$ws.Range("N80").Value = -600005322
$ws.Range("H83").Value = 143143730
$ws.Range("I83").Value = 1000151
$ws.Range("J83").Value = 200001150
$ws.Range("K83").Value = 9001359
$ws.Range("L83").Value = 1800010350
$ws.Range("M83").Value = -8996679
$ws.Range("N83").Value = -1800019710
$ws.Range("H113").Value = 4737.5835
$ws.Range("I113").Value = 8834.333000000001
$ws.Range("J113").Value = 640.8333
$ws.Range("K113").Value = 26502.999
$ws.Range("L113").Value = 1922.4999
$ws.Range("M113").Value = -24332.999
$ws.Range("N113").Value = -6262.4999
$ws.Range("H122").Value = 3157.1
$ws.Range("I122").Value = 434.8889
$ws.Range("J122").Value = 5384.364
$ws.Range("K122").Value = 3914.0001
$ws.Range("L122").Value = 48459.276
$ws.Range("M122").Value = -1464.0001
$ws.Range("N122").Value = -53359.276
$ws.Range("H135").Value = 4203.8184
$ws.Range("I135").Value = 8008.2856
$ws.Range("J135").Value = 1400.5264
$ws.Range("K135").Value = 72074.5704
$ws.Range("L135").Value = 12604.7376
$ws.Range("M135").Value = -69539.5704
$ws.Range("N135").Value = -17674.7376

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 860632.8
$ws.Range("I2").Value = 1338738
$ws.Range("J2").Value = 43.4
$ws.Range("K2").Value = 1338738
$ws.Range("L2").Value = 43.4
$ws.Range("M2").Value = -1338625
$ws.Range("N2").Value = -269.4
$ws.Range("H132").Value = 1536.7261
$ws.Range("I132").Value = 981.92
$ws.Range("J132").Value = 2742.8262
$ws.Range("K132").Value = 2945.76
$ws.Range("L132").Value = 8228.4786
$ws.Range("M132").Value = -415.7599999999998
$ws.Range("N132").Value = -13288.4786

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 937.44446
$ws.Range("J22").Value = 937.4
$ws.Range("L22").Value = 937.4
$ws.Range("N22").Value = -1527.4
$ws.Range("H27").Value = 937.44446
$ws.Range("J27").Value = 937.4
$ws.Range("L27").Value = 937.4
$ws.Range("N27").Value = -1151.4
$ws.Range("H46").Value = 7812.067
$ws.Range("I46").Value = 3909.875
$ws.Range("J46").Value = 12271.714
$ws.Range("K46").Value = 3909.875
$ws.Range("L46").Value = 12271.714
$ws.Range("M46").Value = -3721.875
$ws.Range("N46").Value = -12647.714
$ws.Range("H122").Value = 2501
$ws.Range("I122").Value = 2269.3333
$ws.Range("J122").Value = 2640
$ws.Range("K122").Value = 6807.999899999999
$ws.Range("L122").Value = 7920
$ws.Range("M122").Value = -4357.999899999999
$ws.Range("N122").Value = -12820
$ws.Range("H132").Value = 2225.4717
$ws.Range("I132").Value = 1545.2069
$ws.Range("J132").Value = 3047.4583
$ws.Range("K132").Value = 4635.620699999999
$ws.Range("L132").Value = 9142.374899999999
$ws.Range("M132").Value = -2105.620699999999
$ws.Range("N132").Value = -14202.3749
$ws.Range("H136").Value = 1482.7572
$ws.Range("I136").Value = 1112.4423
$ws.Range("J136").Value = 2552.5557
$ws.Range("K136").Value = 3337.3269
$ws.Range("L136").Value = 7657.6671
$ws.Range("M136").Value = -787.3269
$ws.Range("N136").Value = -12757.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1505.2098
$ws.Range("I132").Value = 1378.8524
$ws.Range("J132").Value = 1890.6
$ws.Range("K132").Value = 4136.5572
$ws.Range("L132").Value = 5671.799999999999
$ws.Range("M132").Value = -1606.5572
$ws.Range("N132").Value = -10731.8
$ws.Range("H136").Value = 781.46875
$ws.Range("I136").Value = 601.90625
$ws.Range("J136").Value = 961.03125
$ws.Range("K136").Value = 1805.71875
$ws.Range("L136").Value = 2883.09375
$ws.Range("M136").Value = 744.28125
$ws.Range("N136").Value = -7983.09375
